$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename sheet
$ws.Name = "sheet1"

# Motor Sizing Inputs
$ws.Range("B4").Value = 25
$ws.Range("B5").Value = 600
$ws.Range("B6").Value = 64
$ws.Range("D6").Value = "64:1 gear reduction"

# Gears
$ws.Range("B14").Value = 42
$ws.Range("B15").Value = 21
$ws.Range("B16").Value = 16
$ws.Range("B17").Value = 42

# Wheel
$ws.Range("B20").Value = 38
$ws.Range("B21").Value = 2.4

# Springs
$ws.Range("B25").Formula = "=(25/2.205*9.81)/(100-20)"
$ws.Range("B26").Value = 100
$ws.Range("B27").Value = 20
$ws.Range("B29").Value = 2

# Shafts
$ws.Range("B32").Value = 145
$ws.Range("B33").Value = 3

# Spring Sizing (formula replaced by a hardcoded constant)
$ws.Range("B37").Value = 80

# Highlight fill for changed/highlighted cells (applied individually so Excel
# reuses a single merged style rather than the comma-list range quirk)
foreach ($addr in @("B14","B16","B17","B20","B32","B33","B37")) {
    $ws.Range($addr).Interior.Color = 65535
}
